# Update "想去人数" (want-to-go count) values in column F for the sheets
# that contain the 漫展 (convention) listing data: "展览" and "全部类型".
# Both sheets mirror the same rows, so the same updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 1555
    6  = 1085
    7  = 11277
    10 = 431
    14 = 12292
    15 = 12927
    22 = 76
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
